# Generate Report for Handoff
# Adds a new localized file entry (a83a2aa5-c998-475b-be6c-7581a4ddc8cf) as
# row 3 on the "Overview", "zh-cn" and "de-de" worksheets, mirroring the
# existing row 2 entry for 2f04207e-aabe-406d-8df0-6c2b57533827.

$wb = $excel.ActiveWorkbook

$newMdName  = "a83a2aa5-c998-475b-be6c-7581a4ddc8cf.md"
$newZhXlf   = "a83a2aa5-c998-475b-be6c-7581a4ddc8cf.975bf7ea07d1e99ea1bb04ee46e13a09a0f3bca4.zh-cn.xlf"
$newDeXlf   = "a83a2aa5-c998-475b-be6c-7581a4ddc8cf.975bf7ea07d1e99ea1bb04ee46e13a09a0f3bca4.de-de.xlf"

$status     = "Ready for handoff"
$include    = "Include"
$epoch      = "0001-01-01 00:00:00"
$mdExt      = ".md"

$overviewDate = "2016-03-24 00:39:12"
$zhHandoffDt  = "2016-03-24 00:39:08"
$deHandoffDt  = "2016-03-24 00:39:12"

$mdBlobUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/46e14cbaf411e4298f4540742e03439875c8cdbc/e2e/$newMdName"
$zhBlobUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d501eaa22b62814a31e591d9b67c81dbfb421452/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newZhXlf"
$deBlobUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9eaa48d370b425b188d22f605b2acf479da12aff/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newDeXlf"

# ---------------------------------------------------------------------
# Overview sheet: File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = $newMdName
$wsOverview.Range("B3").Value = $status
$wsOverview.Range("C3").Value = $status
$wsOverview.Range("D3").Value = $overviewDate
$wsOverview.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $mdBlobUrl, "", "", $newMdName)
$wsOverview.Range("A3").Style = "HyperLink"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A3").Value = $newMdName
$wsZh.Range("B3").Value = $mdExt
$wsZh.Range("C3").Value = $status
$wsZh.Range("D3").Value = $newZhXlf
$wsZh.Range("E3").Value = $zhHandoffDt
$wsZh.Range("H3").Value = $epoch
$wsZh.Range("J3").Value = $include

$wsZh.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $mdBlobUrl, "", "", $newMdName)
$wsZh.Range("A3").Style = "HyperLink"
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), $zhBlobUrl, "", "", $newZhXlf)
$wsZh.Range("D3").Style = "HyperLink"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A3").Value = $newMdName
$wsDe.Range("B3").Value = $mdExt
$wsDe.Range("C3").Value = $status
$wsDe.Range("D3").Value = $newDeXlf
$wsDe.Range("E3").Value = $deHandoffDt
$wsDe.Range("H3").Value = $epoch
$wsDe.Range("J3").Value = $include

$wsDe.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $mdBlobUrl, "", "", $newMdName)
$wsDe.Range("A3").Style = "HyperLink"
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), $deBlobUrl, "", "", $newDeXlf)
$wsDe.Range("D3").Style = "HyperLink"

Write-Output "Handoff report rows added for $newMdName"
